$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column E header date: 01 Jan 2024 (serial 45292) - same display format as C1/D1
$ws.Range("E1").Value = 45292
$ws.Range("E1").NumberFormat = $ws.Range("D1").NumberFormat

# Attendance values for column E (Present / Reason), matching the
# "Present, Absent, Reason" data validation list already on C2:E8
$ws.Range("E2").Value = "Reason"
$ws.Range("E3").Value = "Reason"
$ws.Range("E4").Value = "Present"
$ws.Range("E5").Value = "Reason"
$ws.Range("E6").Value = "Present"
$ws.Range("E7").Value = "Present"
$ws.Range("E8").Value = "Present"

# Comments explaining the "Reason" entries in column E, matching the style
# of the existing "Hp:" comments on D3/D7 ("Hp:" author line + note text)
$ws.Range("E2").AddComment("Hp:" + [char]10 + "Family Emergency") | Out-Null
$ws.Range("E3").AddComment("Hp:" + [char]10 + "Out of Town") | Out-Null
$ws.Range("E5").AddComment("Hp:" + [char]10 + "Outside the Pune") | Out-Null

# Selection moved to G10
$ws.Range("G10").Select()
